# Applies the SimCase5_zsim_SimRun3 data refresh:
#  - re-randomizes column A values for rows 6-201 that changed
#  - truncates the trailing rows 202-251 (data now ends at row 201)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(11, 1).Value = 2
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(17, 1).Value = 2
$ws.Cells.Item(18, 1).Value = 2
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(21, 1).Value = 3
$ws.Cells.Item(24, 1).Value = 2
$ws.Cells.Item(25, 1).Value = 2
$ws.Cells.Item(26, 1).Value = 2
$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(30, 1).Value = 3
$ws.Cells.Item(31, 1).Value = 2
$ws.Cells.Item(32, 1).Value = 3
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(35, 1).Value = 3
$ws.Cells.Item(37, 1).Value = 3
$ws.Cells.Item(39, 1).Value = 1
$ws.Cells.Item(42, 1).Value = 3
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(47, 1).Value = 1
$ws.Cells.Item(49, 1).Value = 2
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(51, 1).Value = 2
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(53, 1).Value = 3
$ws.Cells.Item(54, 1).Value = 2
$ws.Cells.Item(55, 1).Value = 2
$ws.Cells.Item(56, 1).Value = 3
$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(58, 1).Value = 3
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(60, 1).Value = 2
$ws.Cells.Item(61, 1).Value = 3
$ws.Cells.Item(64, 1).Value = 3
$ws.Cells.Item(65, 1).Value = 2
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(68, 1).Value = 2
$ws.Cells.Item(69, 1).Value = 2
$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(71, 1).Value = 1
$ws.Cells.Item(72, 1).Value = 1
$ws.Cells.Item(74, 1).Value = 3
$ws.Cells.Item(75, 1).Value = 2
$ws.Cells.Item(76, 1).Value = 2
$ws.Cells.Item(77, 1).Value = 2
$ws.Cells.Item(78, 1).Value = 2
$ws.Cells.Item(79, 1).Value = 3
$ws.Cells.Item(87, 1).Value = 3
$ws.Cells.Item(88, 1).Value = 3
$ws.Cells.Item(89, 1).Value = 1
$ws.Cells.Item(92, 1).Value = 1
$ws.Cells.Item(93, 1).Value = 2
$ws.Cells.Item(94, 1).Value = 3
$ws.Cells.Item(96, 1).Value = 3
$ws.Cells.Item(98, 1).Value = 1
$ws.Cells.Item(100, 1).Value = 1
$ws.Cells.Item(102, 1).Value = 1
$ws.Cells.Item(104, 1).Value = 3
$ws.Cells.Item(105, 1).Value = 1
$ws.Cells.Item(108, 1).Value = 2
$ws.Cells.Item(109, 1).Value = 1
$ws.Cells.Item(110, 1).Value = 1
$ws.Cells.Item(112, 1).Value = 1
$ws.Cells.Item(113, 1).Value = 1
$ws.Cells.Item(114, 1).Value = 3
$ws.Cells.Item(121, 1).Value = 2
$ws.Cells.Item(122, 1).Value = 3
$ws.Cells.Item(123, 1).Value = 1
$ws.Cells.Item(124, 1).Value = 3
$ws.Cells.Item(126, 1).Value = 2
$ws.Cells.Item(127, 1).Value = 1
$ws.Cells.Item(130, 1).Value = 1
$ws.Cells.Item(132, 1).Value = 3
$ws.Cells.Item(133, 1).Value = 2
$ws.Cells.Item(134, 1).Value = 3
$ws.Cells.Item(135, 1).Value = 3
$ws.Cells.Item(138, 1).Value = 3
$ws.Cells.Item(139, 1).Value = 2
$ws.Cells.Item(142, 1).Value = 1
$ws.Cells.Item(145, 1).Value = 3
$ws.Cells.Item(147, 1).Value = 3
$ws.Cells.Item(149, 1).Value = 2
$ws.Cells.Item(150, 1).Value = 3
$ws.Cells.Item(151, 1).Value = 2
$ws.Cells.Item(152, 1).Value = 1
$ws.Cells.Item(154, 1).Value = 3
$ws.Cells.Item(155, 1).Value = 2
$ws.Cells.Item(156, 1).Value = 3
$ws.Cells.Item(157, 1).Value = 3
$ws.Cells.Item(158, 1).Value = 3
$ws.Cells.Item(160, 1).Value = 1
$ws.Cells.Item(161, 1).Value = 1
$ws.Cells.Item(162, 1).Value = 1
$ws.Cells.Item(163, 1).Value = 2
$ws.Cells.Item(166, 1).Value = 1
$ws.Cells.Item(167, 1).Value = 2
$ws.Cells.Item(168, 1).Value = 3
$ws.Cells.Item(169, 1).Value = 1
$ws.Cells.Item(171, 1).Value = 3
$ws.Cells.Item(172, 1).Value = 2
$ws.Cells.Item(173, 1).Value = 3
$ws.Cells.Item(174, 1).Value = 3
$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(179, 1).Value = 2
$ws.Cells.Item(180, 1).Value = 1
$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(182, 1).Value = 2
$ws.Cells.Item(183, 1).Value = 3
$ws.Cells.Item(184, 1).Value = 1
$ws.Cells.Item(185, 1).Value = 2
$ws.Cells.Item(186, 1).Value = 3
$ws.Cells.Item(187, 1).Value = 3
$ws.Cells.Item(188, 1).Value = 1
$ws.Cells.Item(190, 1).Value = 1
$ws.Cells.Item(191, 1).Value = 3
$ws.Cells.Item(193, 1).Value = 1
$ws.Cells.Item(194, 1).Value = 3
$ws.Cells.Item(196, 1).Value = 1
$ws.Cells.Item(197, 1).Value = 2
$ws.Cells.Item(198, 1).Value = 1
$ws.Cells.Item(199, 1).Value = 3
$ws.Cells.Item(201, 1).Value = 3

# Remove the now-truncated tail of the simulation run
$ws.Range("A202:A251").EntireRow.Delete()
